# Rename the three inline "logo" pictures (two Pearson logos living in the
# default/first-page footers, one BTEC logo living in the first-page header).
#
# InlineShape has no settable .Name in the Word object model, so each
# picture is briefly promoted to a floating Shape (ConvertToShape), renamed
# there (which is what actually drives wp:docPr/@name), and then converted
# back to an inline picture (ConvertToInlineShape) so the drawing stays
# wp:inline exactly as it was.

$d   = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footer (default / primary) : Pearson logo -> image2.png -------------
$footerPrimary = $sec.Footers.Item(1)
$pearson1 = $footerPrimary.Range.InlineShapes.Item(1)
$pearson1Shape = $pearson1.ConvertToShape()
$pearson1Shape.Name = "image2.png"
$pearson1Shape.ConvertToInlineShape() | Out-Null

# --- Footer (first page) : Pearson logo -> image2.png ---------------------
$footerFirst = $sec.Footers.Item(2)
$pearson2 = $footerFirst.Range.InlineShapes.Item(1)
$pearson2Shape = $pearson2.ConvertToShape()
$pearson2Shape.Name = "image2.png"
$pearson2Shape.ConvertToInlineShape() | Out-Null

# --- Header (first page) : BTEC logo -> image1.jpg -------------------------
$headerFirst = $sec.Headers.Item(2)
$btec = $headerFirst.Range.InlineShapes.Item(1)
$btecShape = $btec.ConvertToShape()
$btecShape.Name = "image1.jpg"
$btecShape.ConvertToInlineShape() | Out-Null

Write-Host "Renamed inline picture shapes (2x image2.png, 1x image1.jpg)"
